$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (181-190 in 1-based Excel rows, corresponding to source rows 180-190)
# Columns: A=date, B=city, C=province, D=total_cases (formula), E=new_cases, H=map_name
$dates   = @(43918,43918,43918,43918,43918,43918,43918,43918,43918,43918,43918)
$cities  = @("Riyadh","Jeddah","Mecca","Qatif","Medinah","Tabuk","Khamis Mushait","Abha","Ahsaa","Kobar","Qatif")
$provs   = @("Riyadh","Makkah","Makkah","Eastern province","Medinah","Tabuk","Asir","Asir","Eastern province","Eastern province","Eastern province")
$newCase = @(41,18,12,12,6,3,3,1,1,1,1)
$mapName = @("Ar Riyad","Makkah","Makkah","Ash Sharqiyah","Al Madinah","Tabouk",'`Asir','`Asir',"Ash Sharqiyah","Ash Sharqiyah","Ash Sharqiyah")

$startRow = 180
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd"
    $ws.Cells.Item($r, 2).Value = $cities[$i]
    $ws.Cells.Item($r, 3).Value = $provs[$i]
    $ws.Cells.Item($r, 5).Value = $newCase[$i]
    $ws.Cells.Item($r, 8).Value = $mapName[$i]
}

# Fill column D with the running-total formula, relative per row (matches D136+E137 pattern)
$ws.Range("D180:D190").FormulaR1C1 = "=R[-1]C+RC[1]"

# Update the sheet view to reflect the new extent / selection
$ws.Application.ActiveWindow.ScrollRow = 159
$ws.Range("G192").Select()
